$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the date in J1 (was 43139 / 2018-02-08, now 43161 / 2018-03-02)
$ws.Range("J1").Value2 = 43161

# Add a new row of data (row 15)
$ws.Range("A15").Value = "Dokumentazon "
$ws.Range("I15").Value = "S"
$ws.Range("J15").Value = "S"

# Update the selection to match the author's final cursor position
$ws.Range("I17").Select()
